# Update odds/match data on "Jogos da Semana" sheet to match the
# 2024-11-23 FlashScore refresh (updated odds for several matches plus a
# full data refresh for row 3 - Royal Pari vs Wilstermann, replacing the
# previous Tomayapo vs SA Bulo Bulo entry).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Tigre vs Instituto) - a few odds tweaks
$ws.Range("I2").Value = 3.25
$ws.Range("AD2").Value = 5.5
$ws.Range("AJ2").Value = 12
$ws.Range("BD2").Value = 126

# Row 3 - replaced entirely: new match id/time/teams and a full odds refresh
$ws.Range("A3").Value = "8Ecsv645"
$ws.Range("C3").Value = "21:00"
$ws.Range("E3").Value = "Royal Pari"
$ws.Range("F3").Value = "Wilstermann"
$ws.Range("G3").Value = 2.3
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 3.75
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.67
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 11
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 21
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 8.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 17
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 351
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 15
$ws.Range("AJ3").Value = 11
$ws.Range("AK3").Value = 34
$ws.Range("AM3").Value = 41
$ws.Range("AN3").Value = 4.33
$ws.Range("AO3").Value = 13
$ws.Range("AP3").Value = 26
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 67
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.5
$ws.Range("AU3").Value = 8.5
$ws.Range("AW3").Value = 5
$ws.Range("AY3").Value = 29
$ws.Range("BB3").Value = 251

# Row 4 (Atletico GO vs Palmeiras) - a few odds tweaks
$ws.Range("M4").Value = 1.03
$ws.Range("O4").Value = 1.19
$ws.Range("Q4").Value = 1.73
$ws.Range("R4").Value = 2.08

# Row 5 (Botafogo RJ vs Vitoria) - odds tweaks
$ws.Range("G5").Value = 1.36
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 1.83
$ws.Range("K5").Value = 2.4
$ws.Range("L5").Value = 8
$ws.Range("Y5").Value = 9
$ws.Range("AD5").Value = 9.5
$ws.Range("AE5").Value = 23
$ws.Range("AF5").Value = 81
$ws.Range("AG5").Value = 501
$ws.Range("AH5").Value = 19
$ws.Range("AL5").Value = 67
$ws.Range("AN5").Value = 3.2
$ws.Range("AU5").Value = 10

# Row 6 (Juventude vs Cuiaba) - odds tweaks
$ws.Range("G6").Value = 1.62
$ws.Range("I6").Value = 6
$ws.Range("K6").Value = 2.1
$ws.Range("L6").Value = 6.5
$ws.Range("R6").Value = 1.65
$ws.Range("Z6").Value = 11
$ws.Range("AC6").Value = 8
$ws.Range("AO6").Value = 8.5
$ws.Range("AZ6").Value = 151
$ws.Range("BA6").Value = 201

# Row 7 (New York City vs New York Red Bulls) - odds tweaks
$ws.Range("G7").Value = 2.1
$ws.Range("I7").Value = 3.4
$ws.Range("J7").Value = 2.75
$ws.Range("L7").Value = 4
$ws.Range("O7").Value = 1.3
$ws.Range("P7").Value = 3.4
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.85
$ws.Range("U7").Value = 1.75
$ws.Range("V7").Value = 2
$ws.Range("W7").Value = 8
$ws.Range("X7").Value = 10
$ws.Range("Z7").Value = 19
$ws.Range("AC7").Value = 10
$ws.Range("AE7").Value = 15
$ws.Range("AF7").Value = 51
$ws.Range("AG7").Value = 251
$ws.Range("AJ7").Value = 12
$ws.Range("AK7").Value = 41
$ws.Range("AN7").Value = 4
$ws.Range("AO7").Value = 11
$ws.Range("AU7").Value = 8
$ws.Range("AW7").Value = 5.5
$ws.Range("AX7").Value = 19
$ws.Range("AZ7").Value = 67
$ws.Range("BB7").Value = 201
